$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Rename existing headers I1/J1 (pandas merge suffix artifact)
$ws.Cells.Item(1, 9).Value = "DS_ESTADO_CIVIL_x"
$ws.Cells.Item(1, 10).Value = "DS_GRAU_INSTRUCAO_x"

# New header cells K1, L1, M1, N1 - copy the formatting from I1 (bold/border/center style)
# then overwrite with the correct header text.
$ws.Cells.Item(1, 9).Copy($ws.Cells.Item(1, 11))
$ws.Cells.Item(1, 9).Copy($ws.Cells.Item(1, 12))
$ws.Cells.Item(1, 9).Copy($ws.Cells.Item(1, 13))
$ws.Cells.Item(1, 9).Copy($ws.Cells.Item(1, 14))

$ws.Cells.Item(1, 11).Value = "DS_ESTADO_CIVIL_y"
$ws.Cells.Item(1, 12).Value = "DS_GRAU_INSTRUCAO_y"
$ws.Cells.Item(1, 13).Value = "SQ_CANDIDATO"
$ws.Cells.Item(1, 14).Value = "SG_PARTIDO"

# --- Data rows (2-85) ---
# K = copy of I (DS_ESTADO_CIVIL), L = copy of J (DS_GRAU_INSTRUCAO), M = copy of F (SQ_CANDIDATO_OF)
$ws.Range("I2:I85").Copy($ws.Range("K2:K85"))
$ws.Range("J2:J85").Copy($ws.Range("L2:L85"))
$ws.Range("F2:F85").Copy($ws.Range("M2:M85"))

# N = SG_PARTIDO (new data, one value per candidate row)
$parties = @(
    "PSL",
    "PSL",
    "PSL",
    "PSL",
    "PATRIOTA",
    "PSL",
    "PSL",
    "PSL",
    "PSL",
    "PSC",
    "PROS",
    "PROS",
    "PRP",
    "PSOL",
    "PSOL",
    "PSOL",
    "PSOL",
    "PV",
    "PV",
    "PRB",
    "PODE",
    "PODE",
    "PRP",
    "PMB",
    "PPL",
    "REDE",
    "REDE",
    "PRP",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PDT",
    "PP",
    "DEM",
    "PR",
    "PTB",
    "PSDB",
    "PMN",
    "MDB",
    "MDB",
    "MDB",
    "PP",
    "PTC",
    "PSD",
    "MDB",
    "PPS",
    "PP",
    "SOLIDARIEDADE",
    "PT",
    "PT",
    "PT",
    "PT",
    "PSB",
    "PT",
    "PC do B",
    "PHS",
    "PHS",
    "DC",
    "DC",
    "DC",
    "AVANTE",
    "AVANTE",
    "AVANTE",
    "AVANTE",
    "AVANTE",
    "AVANTE",
    "AVANTE",
    "AVANTE",
    "AVANTE",
    "PRTB",
    "PRTB",
    "PDT"
)
for ($i = 0; $i -lt $parties.Length; $i++) {
    $ws.Cells.Item($i + 2, 14).Value = $parties[$i]
}

Write-Host "Done: added columns K-N (rows 1-85)"
